$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attributes Sheet")
$ws.Activate()

# --- Row 2 ("تقرير 39" / Document Title) ---
# Mandatory Y -> N
$ws.Range("C2").Value = "N"
# Max Length 244 -> 255
$ws.Range("E2").Value = 255
# Default value: (blank) -> Test
$ws.Range("J2").Value = "Test"

# --- Row 3 ("تقرير 39" / حاله المستند) ---
# Mandatory N -> Y
$ws.Range("C3").Value = "Y"

# --- Default value "Test" added to the other "Document Title" rows ---
$ws.Range("J6").Value = "Test"
$ws.Range("J11").Value = "Test"
$ws.Range("J16").Value = "Test"
$ws.Range("J21").Value = "Test"
$ws.Range("J24").Value = "Test"
$ws.Range("J27").Value = "Test"

# --- Row 7 (Oracle حوافظ / حاله المستند): Default value (blank) -> Not Approved ---
$ws.Range("J7").Value = "Not Approved"

# --- Update the sheet view: scroll so row 7 is at the top, select J12 ---
$ws.Range("J12").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 2
